# Updates the Price (D) and Volume(1h) (E) columns for this crypto-price
# snapshot, matching the refreshed values from the scheduled GitHub Actions
# scrape run on Thu Jan 19 14:48:32 UTC 2023.
#
# Values are written with a leading apostrophe so Excel keeps them as text
# (these columns store formatted strings like "289.13" / "-4.58%", not
# numeric types), matching how the sheet already stores this data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'289.13"
$ws.Range("E2").Value = "'-4.58%"
$ws.Range("D3").Value = "'30.23"
$ws.Range("E3").Value = "'-7.37%"
$ws.Range("D4").Value = "'4.924"
$ws.Range("E4").Value = "'-3.03%"
$ws.Range("D5").Value = "'0.07244"
$ws.Range("E5").Value = "'-6.14%"
$ws.Range("D6").Value = "'1.798"
$ws.Range("E6").Value = "'-13.79%"
$ws.Range("D7").Value = "'7.591"
$ws.Range("E7").Value = "'-3.97%"
$ws.Range("D8").Value = "'3.706"
$ws.Range("E8").Value = "'-2.44%"
$ws.Range("D9").Value = "'0.9035"
$ws.Range("E9").Value = "'-2.10%"
$ws.Range("D10").Value = "'0.1676"
$ws.Range("E10").Value = "'-5.14%"
$ws.Range("D11").Value = "'0.07845"
$ws.Range("E11").Value = "'-1.67%"
$ws.Range("D12").Value = "'0.07995"
$ws.Range("E12").Value = "'-6.95%"
$ws.Range("D13").Value = "'0.03045"
$ws.Range("E13").Value = "'-0.54%"
$ws.Range("D14").Value = "'0.1002"
$ws.Range("E14").Value = "'0.33%"
$ws.Range("D15").Value = "'0.001490"
$ws.Range("E15").Value = "'-1.73%"
$ws.Range("D16").Value = "'0.005740"
$ws.Range("E16").Value = "'1.15%"
$ws.Range("D17").Value = "'3.467"
$ws.Range("E17").Value = "'-0.08%"
$ws.Range("D18").Value = "'2.075"
$ws.Range("E18").Value = "'-3.78%"
$ws.Range("E19").Value = "'-0.97%"
$ws.Range("E20").Value = "'-0.52%"
$ws.Range("D21").Value = "'3.974"
$ws.Range("E21").Value = "'-9.55%"
$ws.Range("E22").Value = "'10.03%"
$ws.Range("D23").Value = "'0.04503"
$ws.Range("E23").Value = "'-0.71%"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'-1.53%"
$ws.Range("D25").Value = "'0.004436"
$ws.Range("E25").Value = "'7.04%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'3.86%"
$ws.Range("D39").Value = "'0.01567"
$ws.Range("E39").Value = "'-9.19%"
$ws.Range("D40").Value = "'0.04346"
$ws.Range("D41").Value = "'0.007280"
$ws.Range("E41").Value = "'-2.70%"
$ws.Range("D42").Value = "'0.01007"
$ws.Range("D43").Value = "'0.1305"
$ws.Range("E43").Value = "'-4.07%"
$ws.Range("E44").Value = "'-14.07%"
$ws.Range("D45").Value = "'0.009040"
$ws.Range("E45").Value = "'-14.76%"
$ws.Range("D46").Value = "'0.00005889"
$ws.Range("E46").Value = "'-5.12%"
$ws.Range("E47").Value = "'-0.18%"
$ws.Range("E48").Value = "'104.79%"
$ws.Range("E49").Value = "'-3.62%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.18%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.18%"
Write-Output "Updated 39 Price cells and 49 Volume(1h) cells on $($ws.Name)."
